$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: a duplicate of row 2's order/product data, marked as having
# "Pajareras" (bird netting) = "Si" this time.
$ws.Range("A3").Value = $ws.Range("A2").Value()
$ws.Range("B3").Value = $ws.Range("B2").Value()
$ws.Range("C3").Value = $ws.Range("C2").Value()
$ws.Range("D3").Value = "'" + $ws.Range("D2").Value()
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = ""
$ws.Range("H3").Value = $ws.Range("H2").Value()
$ws.Range("I3").Value = "'" + $ws.Range("I2").Value()
$ws.Range("J3").Value = $ws.Range("J2").Value()
$ws.Range("K3").Value = "'" + $ws.Range("K2").Value()
$ws.Range("L3").Value = ""
$ws.Range("M3").Value = "Sí"
$ws.Range("N3").Value = $ws.Range("N2").Value()

# Row 2's placeholder blanks aren't needed any more now that row 3 exists
$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()
$ws.Range("G2").ClearContents()
$ws.Range("L2").ClearContents()
$ws.Range("M2").ClearContents()
